# coppel-data.xlsx update:
#  - "agregando resgistros al excel"
#  - MCP-INT-01 sheet: row 3 H column becomes an email (hyperlink) instead of a
#    phone number, and a new row 13 with G13="S" is added.
#  - MCP-INT-02 sheet: a new data row (row 3) is appended, and it becomes the
#    active/selected sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("MCP-INT-01")
$ws2 = $wb.Worksheets.Item("MCP-INT-02")

# --- MCP-INT-01 ("Sheet2" in the package) -----------------------------------

# H3 used to hold a phone number; it now holds an e-mail address that is also
# turned into a mailto: hyperlink (mirroring the existing H2 hyperlink), with
# the same visual style (Hyperlink + vertical-top alignment + unlocked).
$h3 = $ws1.Range("H3")
$h3.Value = "hrpinam@gmail.com"
$h3.Style = "Hyperlink"
$ws1.Hyperlinks.Add($h3, "mailto:hrpinam@gmail.com")
$h3.VerticalAlignment = -4160
$h3.Locked = $false

# A brand new row with a single populated cell.
$ws1.Range("G13").Value = "S"

# Update the selection on MCP-INT-01 before moving away from it.
[void]$ws1.Activate()
[void]$ws1.Range("C26").Select()

# --- MCP-INT-02 ("Sheet3" in the package) -----------------------------------

$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = "MCP-INT-02"
$ws2.Range("C3").Value = "Chrome"
$ws2.Range("D3").Value = "https://www.coppel.com/"
$ws2.Range("E3").Value = "Maria"
$ws2.Range("F3").Value = "Test"
$ws2.Range("G3").Value = "Femenino"
$ws2.Range("H3").Value = 6672529402
$ws2.Range("I3").Value = "contra123contra"
$ws2.Range("J3").Value = "Si"

# MCP-INT-02 becomes the active/selected sheet.
[void]$ws2.Activate()
[void]$ws2.Range("J10").Select()
